# "fixed projection, 180 overlap issue remaining"
#
# The sign convention for the track-width delta vector (B1/B2 on the
# "Track limit points" sheet) was flipped, which was causing the
# projected points to land on the wrong side of the track centre line.
# Swap B1 and B2 so the delta-x / delta-y pair point the other way, and
# leave the cursor on the cell the author was inspecting next.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Track limit points")
$ws.Select()

$ws.Range("B1").Value = -2
$ws.Range("B2").Value = 2

$ws.Range("G14").Select()
